# "fix contracts: use names from domain model"
#
# The contract texts for the two system operations (getLocalWorkbooks and
# getLocalPages) are updated so the return/element types match the names
# used in the domain model: "Arbeitsheft" (instead of no return type) and
# "Seite" (instead of the English "Page").

$d = $word.ActiveDocument

# 1) getLocalWorkbooks() contract: add the missing return type.
$d.Content.Find.Execute(
    "getLocalWorkbooks()", $false, $false, $false, $false, $false,
    $true, 1, $false, "getLocalWorkbooks() : Arbeitsheft[]", 2)

# 2) getLocalPages(...) contract: rename the returned type from the
#    English "Page" to the domain model's "Seite".
$d.Content.Find.Execute(
    "getLocalPages(workbook : Workbook) : Page[]", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "getLocalPages(workbook : Workbook) : Seite[]", 2)

# 3) Matching "Ergebnisse" (results) description for getLocalPages: use
#    "Seiten" instead of "Page-Objekte".
$d.Content.Find.Execute(
    "Die Page-Objekte, die mit workbook assoziiert sind", $false, $false,
    $false, $false, $false, $true, 1, $false,
    "Die Seiten, die mit workbook assoziiert sind", 2)
